$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the SamplesTab query text (cell B4) with the revised SQL that now
# also selects participant_age_at_collection, sample_tumor_status and
# tumor_classification, and surfaces them as extra output columns.
$newQuery = @'
WITH filtered_participants AS (
    SELECT DISTINCT prt.id, prt.participant_id
    FROM df_study std
    JOIN df_participant prt ON std.id = prt."study.id"
    JOIN df_diagnosis dgn ON prt.id = dgn."participant.id"
    WHERE 
        std.dbgap_accession = 'phs002371'
        AND prt.sex_at_birth = 'Female'
        AND dgn.diagnosis = '9861/3 : Acute myeloid leukemia, NOS'
),
-- Step 2: Get samples from those participants
filtered_samples AS (
    SELECT DISTINCT smp.id, smp.sample_id, smp.anatomic_site, smp."participant.id",smp.participant_age_at_collection, smp.sample_tumor_status,smp.tumor_classification
    FROM df_sample smp
    JOIN filtered_participants prt ON smp."participant.id" = prt.id
),
-- Step 3: Keep only samples that have at least one FASTQ file
samples_with_fastq AS (
    SELECT DISTINCT smp.*
    FROM filtered_samples smp
    JOIN df_sequencing_file seq ON seq."sample.id" = smp.id
    WHERE seq.file_type = 'fastq'
)
-- Final output
SELECT 
    smp.sample_id AS "Sample ID",
    prt.participant_id AS "Participant ID",
    'phs002371' AS "Study ID",
    smp.anatomic_site AS "Sample Anatomic Site",
    COALESCE(CASE WHEN smp.participant_age_at_collection = -999 THEN 'Not Reported' ELSE smp.participant_age_at_collection END, 0) AS "Age at Sample Collection (days)",
    COALESCE(smp.sample_tumor_status, '') AS "Sample Tumor Status",
    COALESCE(smp.tumor_classification, '') AS "Sample Tumor Classification",
    NULL AS "Sample Diagnosis"
FROM 
    samples_with_fastq smp
JOIN filtered_participants prt ON smp."participant.id" = prt.id
ORDER BY 
    smp.sample_id ASC;
'@

$ws.Range("B4").Value = $newQuery

# Move the active view/selection down one row, matching the saved view state.
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
